# Actualización automática de tasas-transfi.xlsx
$wb = $excel.ActiveWorkbook

# --- Update "Hoja1" A1 text block with new conversion rates ---
$hoja1 = $wb.Worksheets.Item("Hoja1")
$cell = $hoja1.Range("A1")
$text = [string]$cell.Value()
$text = $text.Replace("✅ 1000 Bs = 6.12 = 24770.49 pesos", "✅ 1000 Bs = 6.04 = 24776.41 pesos")
$text = $text.Replace("✅ 24770.49 pesos = 6.08 = 962.24 Bs", "✅ 24776.41 pesos = 6.07 = 970.43 Bs")
$cell.Value = $text

# --- Update "tasas" sheet rate values ---
$tasas = $wb.Worksheets.Item("tasas")
$tasas.Range("N10").Value = 165.48
$tasas.Range("O10").Value = 4100
$tasas.Range("N12").Value = 4085
$tasas.Range("O12").Value = 160
